$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data change: state column (F) value corrected from lowercase to
#     proper-case "California" for all data rows ---
$ws.Range("F2:F5").Value = "California"

# --- Column F (new "state" column) sized like the other data columns ---
$ws.Columns.Item(6).ColumnWidth = 20.59

# --- Header cell E1 ("city") keeps its bold weight but loses the green
#     header fill (distinguishing it from the other header cells) ---
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").Interior.ColorIndex = -4142

# --- Data cells E2:E5 ("sunny") explicitly cleared of any fill ---
$ws.Range("E2:E5").Interior.ColorIndex = -4142

# --- View state: scroll the sheet down a bit and move the selection onto
#     the new state column ---
$ws.Activate() | Out-Null
$w = $excel.ActiveWindow
$w.ScrollRow = 4
$w.ScrollColumn = 1
$ws.Range("F2:F5").Select() | Out-Null
